$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")
Write-Host "Sheet name: $($ws.Name)"
